$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update MACRO_SCORE column (N) for rows 2-6
$ws.Range("N2").Value = 85.77505782882612
$ws.Range("N3").Value = 85.77505782882612
$ws.Range("N4").Value = 85.77505782882612
$ws.Range("N5").Value = 85.77505782882612
$ws.Range("N6").Value = 85.77505782882612

# Update 종가 (close price) and 5일수익률 (5-day return) for row 4
$ws.Range("D4").Value = 90922.84
$ws.Range("F4").Value = 4.74
